$d = $word.ActiveDocument

# Phase 1: move every changed [[PERSON_<old>]] to a unique temporary token
# so that cascading renumbering does not clobber already-written values.
$d.Content.Find.Execute("[[PERSON_11]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_11]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_12]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_12]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_13]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_13]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_14]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_14]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_15]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_15]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_16]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_16]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_17]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_17]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_18]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_18]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_19]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_19]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_20]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_20]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_21]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_21]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_22]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_22]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_23]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_23]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_24]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_24]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_25]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_25]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_26]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_26]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_27]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_27]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_28]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_28]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_29]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_29]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_30]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_30]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_31]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_31]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_32]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_32]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_33]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_33]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_34]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_34]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_35]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_35]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_36]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_36]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_37]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_37]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_38]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_38]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_39]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_39]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_40]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_40]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_41]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_41]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_42]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_42]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_43]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_43]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_44]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_44]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_45]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_45]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_46]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_TMP_46]]", 2) | Out-Null

# Phase 2: move every temporary token to its final new number
$d.Content.Find.Execute("[[PERSON_TMP_11]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_10]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_12]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_11]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_13]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_12]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_14]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_13]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_15]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_14]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_16]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_15]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_17]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_16]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_18]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_17]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_19]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_18]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_20]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_19]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_21]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_20]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_22]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_21]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_23]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_22]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_24]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_23]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_25]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_24]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_26]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_25]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_27]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_26]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_28]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_27]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_29]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_28]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_30]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_29]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_31]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_30]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_32]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_30]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_33]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_31]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_34]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_32]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_35]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_33]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_36]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_34]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_37]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_35]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_38]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_36]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_39]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_37]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_40]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_38]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_41]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_39]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_42]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_40]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_43]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_40]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_44]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_41]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_45]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_42]]", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_TMP_46]]", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_26]]", 2) | Out-Null
